# Insert a new data row at row 776 (pushing the existing rows 776:849 down
# to 777:850) and populate it with a new price observation.
#
# Copying row 776 first and inserting the copy means the new row inherits
# the same formatting/number-format/style (e.g. the date style on column D)
# as the surrounding data rows, and all of the fields that don't change
# (Mercado ID, Mercado, Region, Codreg, Categoria ID, Categoria, Variedad,
# Calidad, Kg o Unidades, Clasificacion) are already correct because they
# are identical to the row being copied.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(776).Copy()
$ws.Rows(776).Insert()

$ws.Range("D776").Value = 45194
$ws.Range("J776").Value = 34
$ws.Range("K776").Value = 15000
$ws.Range("L776").Value = 16000
$ws.Range("M776").Value = 15500
$ws.Range("N776").Value = "$/malla 25 kilos"
$ws.Range("O776").Value = "Perú"
$ws.Range("P776").Value = 620
